$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the runNumber column (D) values from 4717 to 4715 for rows 2 through 63
$ws.Range("D2:D63").Value = 4715

# Reflect the resulting selection on the sheet (D3:D63, active cell D3)
$ws.Range("D3:D63").Select()
